$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column F - reuse the existing bold/bordered header style
# (copy format only, so no new font/style entries get created in styles.xml)
$ws.Cells.Item(1, 6).Value = "Trening"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Rewrite column A (Timestamp) from text to real datetimes, refresh all data rows,
# and add new column F ("Trening" - training part / Duza Gra vs Mala Gra),
# plus append the new rows captured for the second half of the session.
$ws.Cells.Item(2, 1).Value = 45685.64817534723
$ws.Cells.Item(2, 2).Value = 1117.3
$ws.Cells.Item(2, 3).Value = 14.3
$ws.Cells.Item(2, 4).Value = 3.912192991801672
$ws.Cells.Item(2, 5).Value = "10-15"
$ws.Cells.Item(2, 6).Value = "Duża Gra"
$ws.Cells.Item(3, 1).Value = 45685.64953182871
$ws.Cells.Item(3, 2).Value = 1234.5
$ws.Cells.Item(3, 3).Value = 14.08
$ws.Cells.Item(3, 4).Value = 4.290932144437516
$ws.Cells.Item(3, 5).Value = "10-15"
$ws.Cells.Item(3, 6).Value = "Duża Gra"
$ws.Cells.Item(4, 1).Value = 45685.65100636574
$ws.Cells.Item(4, 2).Value = 1361.9
$ws.Cells.Item(4, 3).Value = 14.35
$ws.Cells.Item(4, 4).Value = 4.229175295148576
$ws.Cells.Item(4, 5).Value = "10-15"
$ws.Cells.Item(4, 6).Value = "Duża Gra"
$ws.Cells.Item(5, 1).Value = 45685.648171875
$ws.Cells.Item(5, 2).Value = 1117
$ws.Cells.Item(5, 3).Value = 9.74
$ws.Cells.Item(5, 4).Value = 3.358905383518764
$ws.Cells.Item(5, 5).Value = "5-10"
$ws.Cells.Item(5, 6).Value = "Duża Gra"
$ws.Cells.Item(6, 1).Value = 45685.64952835648
$ws.Cells.Item(6, 2).Value = 1234.2
$ws.Cells.Item(6, 3).Value = 8.85
$ws.Cells.Item(6, 4).Value = 3.501034532274519
$ws.Cells.Item(6, 5).Value = "5-10"
$ws.Cells.Item(6, 6).Value = "Duża Gra"
$ws.Cells.Item(7, 1).Value = 45685.65100289352
$ws.Cells.Item(7, 2).Value = 1361.6
$ws.Cells.Item(7, 3).Value = 9.21
$ws.Cells.Item(7, 4).Value = 3.569845982960293
$ws.Cells.Item(7, 5).Value = "5-10"
$ws.Cells.Item(7, 6).Value = "Duża Gra"
$ws.Cells.Item(8, 1).Value = 45685.67112442129
$ws.Cells.Item(8, 2).Value = 3100.1
$ws.Cells.Item(8, 3).Value = 12.12
$ws.Cells.Item(8, 4).Value = 3.78263885634286
$ws.Cells.Item(8, 5).Value = "10-15"
$ws.Cells.Item(8, 6).Value = "Mała Gra"
$ws.Cells.Item(9, 1).Value = 45685.6753130787
$ws.Cells.Item(9, 2).Value = 3462
$ws.Cells.Item(9, 3).Value = 14.41
$ws.Cells.Item(9, 4).Value = 3.600886072431294
$ws.Cells.Item(9, 5).Value = "10-15"
$ws.Cells.Item(9, 6).Value = "Mała Gra"
$ws.Cells.Item(10, 1).Value = 45685.67829108796
$ws.Cells.Item(10, 2).Value = 3719.3
$ws.Cells.Item(10, 3).Value = 14.29
$ws.Cells.Item(10, 4).Value = 3.487117052078247
$ws.Cells.Item(10, 5).Value = "10-15"
$ws.Cells.Item(10, 6).Value = "Mała Gra"
$ws.Cells.Item(11, 1).Value = 45685.66794502315
$ws.Cells.Item(11, 2).Value = 2825.4
$ws.Cells.Item(11, 3).Value = 9.13
$ws.Cells.Item(11, 4).Value = 2.946751492364068
$ws.Cells.Item(11, 5).Value = "5-10"
$ws.Cells.Item(11, 6).Value = "Mała Gra"
$ws.Cells.Item(12, 1).Value = 45685.67112210648
$ws.Cells.Item(12, 2).Value = 3099.9
$ws.Cells.Item(12, 3).Value = 8.83
$ws.Cells.Item(12, 4).Value = 3.442386286599296
$ws.Cells.Item(12, 5).Value = "5-10"
$ws.Cells.Item(12, 6).Value = "Mała Gra"
$ws.Cells.Item(13, 1).Value = 45685.683890625
$ws.Cells.Item(13, 2).Value = 4203.1
$ws.Cells.Item(13, 3).Value = 8.3
$ws.Cells.Item(13, 4).Value = 2.841054993016379
$ws.Cells.Item(13, 5).Value = "5-10"
$ws.Cells.Item(13, 6).Value = "Mała Gra"

# Apply the date/time display format used for the Timestamp column.
# A throwaway lowercase pass registers numFmtId 164 (matches the source file),
# then the real uppercase format is applied across the whole A2:A13 range so
# every timestamp cell shares a single style entry (numFmtId 165).
$ws.Range("A2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("A2:A13").NumberFormat = "YYYY-MM-DD HH:MM:SS"
